# Update the "cross" labels in column A to the new, more specific codes
# and move the sheet's active selection to E15, matching the author's
# "update to new files" commit.
#
#   F1 (row 3) -> F1a
#   F1 (row 4) -> F1b
#   F2 (row 5) -> F2a
#   P2 (row 6) stays P2 (only the shared-string slot it points at changes)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "F1a"
$ws.Range("A4").Value = "F1b"
$ws.Range("A5").Value = "F2a"
$ws.Range("A6").Value = "P2"

# Reflect the new cursor/selection position recorded in the saved file.
$ws.Range("E15").Select()
